$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "maa://24702 (94.54), maa://25390 (96.2), maa://36681 (87.34)"
$ws.Range("D3").Value = "maa://40192 (98.18), maa://36987 (96.08), maa://39849 (88.89)"
$ws.Range("H3").Value = "maa://21247 (98.6), *maa://22748 (60.0)"
$ws.Range("L3").Value = "*maa://22880 (65.33), maa://20276 (86.59), *maa://22749 (75.0)"
$ws.Range("P3").Value = "maa://21249 (94.54), maa://26254 (96.67)"
$ws.Range("T3").Value = "maa://24617 (89.83), **maa://20790 (43.48), ***maa://37170 (16.67), maa://45854 (87.5)"
$ws.Range("X3").Value = "maa://27396 (84.36), maa://27484 (96.58), maa://27480 (83.33)"
$ws.Range("T4").Value = "maa://32509 (95.76), maa://27295 (86.11), maa://22754 (90.41), *maa://21746 (55.81), *maa://31008 (78.57)"
$ws.Range("X4").Value = "**maa://32495 (48.7), ***maa://31785 (22.22), maa://43217 (90.62), ***maa://36683 (28.26)"
$ws.Range("D5").Value = "maa://21245 (84.58), maa://22744 (84.0)"
$ws.Range("L5").Value = "*maa://22757 (76.32)"
$ws.Range("P6").Value = "maa://31836 (93.1), maa://30381 (93.33)"
$ws.Range("P7").Value = "maa://22750 (91.84)"
$ws.Range("X7").Value = "maa://22399 (95.51), *maa://22758 (74.29)"
$ws.Range("A8").Value = "更新日期：2025.03.06 13:18:49"
$ws.Range("P8").Value = "maa://32931 (84.03), *maa://21916 (62.12), maa://23252 (91.18), maa://37496 (96.97), **maa://22759 (45.45)"
$ws.Range("P9").Value = "maa://22736 (83.33)"
$ws.Range("X9").Value = "maa://26223 (97.93)"
$ws.Range("AB9").Value = "maa://28711 (87.1), ***maa://22740 (5.66), **maa://39938 (46.67), **maa://27377 (42.86), ***maa://25174 (19.05), maa://40166 (96.3), *maa://45044 (66.67)"
$ws.Range("AF9").Value = "maa://26206 (87.8), *maa://22865 (51.85)"
$ws.Range("D10").Value = "***maa://25695 (18.62), ***maa://39951 (14.29), ***maa://34206 (19.23), ***maa://39243 (25.0), *maa://45271 (57.89)"
$ws.Range("P10").Value = "maa://28977 (89.41), maa://36669 (86.67), *maa://23264 (61.82)"
$ws.Range("T10").Value = "maa://27395 (96.5), maa://22755 (87.83), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("L11").Value = "maa://21287 (89.42)"
$ws.Range("P11").Value = "maa://45557 (88.89)"
$ws.Range("T11").Value = "maa://22747 (92.5), maa://22501 (97.78), maa://45521 (85.0)"
$ws.Range("H12").Value = "maa://21867 (89.94), **maa://45826 (33.33)"
$ws.Range("X12").Value = "maa://22753 (91.06), *maa://21485 (75.35), maa://37962 (90.48)"
$ws.Range("D13").Value = "maa://24999 (92.16), maa://36673 (93.33), maa://25001 (85.71)"
$ws.Range("H13").Value = "*maa://21248 (73.33), **maa://22728 (46.67)"
$ws.Range("P13").Value = "maa://22676 (92.86), *maa://22583 (75.0), *maa://22500 (58.7)"
$ws.Range("D14").Value = "maa://30764 (87.5)"
$ws.Range("L14").Value = "maa://26245 (96.75), maa://21288 (96.3), maa://39841 (94.55), maa://36682 (97.44)"
$ws.Range("D15").Value = "*maa://22743 (78.04), maa://22734 (84.17), *maa://30808 (64.18), **maa://36048 (46.77), maa://45058 (92.86)"
$ws.Range("H15").Value = "maa://24304 (88.13), maa://21478 (89.19)"
$ws.Range("X15").Value = "*maa://38786 (66.67)"
$ws.Range("AF15").Value = "maa://21364 (80.84), *maa://36666 (77.68), *maa://22766 (68.64)"
$ws.Range("T16").Value = "maa://22729 (94.48), *maa://28648 (69.12), maa://36674 (82.69)"
$ws.Range("D17").Value = "maa://21624 (84.62)"
$ws.Range("H17").Value = "maa://22430 (88.83), maa://39599 (86.27)"
$ws.Range("H18").Value = "maa://24421 (89.02)"
$ws.Range("L18").Value = "maa://22466 (90.53), *maa://22732 (50.55)"
$ws.Range("D20").Value = "maa://21432 (90.23), maa://25198 (93.64), *maa://20795 (50.77), maa://36680 (91.18)"
$ws.Range("H20").Value = "maa://22864 (89.94)"
$ws.Range("L20").Value = "maa://41331 (85.26)"
$ws.Range("L21").Value = "maa://31731 (96.15)"
$ws.Range("AF21").Value = "maa://22524 (94.17), *maa://22432 (77.63)"
$ws.Range("H22").Value = "maa://25236 (95.7), **maa://21678 (48.94), **maa://22735 (42.86)"
$ws.Range("L22").Value = "*maa://27127 (78.99), *maa://22751 (71.01)"
$ws.Range("L23").Value = "maa://39756 (95.74), maa://39875 (94.44)"
$ws.Range("P23").Value = "maa://30587 (92.0), *maa://29748 (75.97), ***maa://29785 (16.18), *maa://37566 (76.92)"
$ws.Range("D24").Value = "*maa://24368 (78.27), **maa://46650 (50.0)"
$ws.Range("X24").Value = "maa://29988 (84.44), maa://23504 (93.21), **maa://22892 (40.14), *maa://25141 (77.1), *maa://36663 (76.92), ***maa://22815 (23.08)"
$ws.Range("AF24").Value = "maa://22523 (85.57), maa://36672 (80.36), maa://29910 (93.22), **maa://21440 (35.71), *maa://45831 (80.0)"
$ws.Range("H25").Value = "*maa://29063 (73.33), *maa://25311 (74.29), ***maa://22725 (4.84), *maa://45047 (62.5)"
$ws.Range("T25").Value = "maa://20109 (92.18), maa://22545 (100.0), *maa://42915 (75.0)"
$ws.Range("AF25").Value = "maa://20108 (96.35), maa://24621 (96.9), maa://36676 (97.06), maa://22771 (85.71), *maa://37772 (66.67)"
$ws.Range("AB26").Value = "maa://42235 (94.23)"
$ws.Range("T28").Value = "*maa://29765 (64.29), maa://23263 (95.28)"
$ws.Range("X28").Value = "maa://39929 (90.63), maa://41749 (90.91), ***maa://39723 (13.89)"
$ws.Range("AF28").Value = "maa://36660 (92.35), *maa://36701 (66.67)"
$ws.Range("L29").Value = "maa://28432 (93.53), maa://28440 (80.18), maa://31400 (98.81), *maa://28650 (71.43)"
$ws.Range("P30").Value = "maa://21442 (99.12)"
$ws.Range("L31").Value = "maa://35926 (93.52), maa://36258 (85.47), *maa://43904 (72.73)"
$ws.Range("H32").Value = "maa://21895 (97.04), maa://36667 (97.62), **maa://20793 (38.78), maa://22760 (100.0)"
$ws.Range("T32").Value = "maa://42859 (95.97), maa://41108 (88.0), maa://41238 (97.14), maa://45523 (100.0)"
$ws.Range("P33").Value = "maa://21956 (81.21), *maa://22730 (79.31)"
$ws.Range("T34").Value = "maa://24526 (93.1)"
$ws.Range("L35").Value = "maa://41296 (96.3)"
$ws.Range("AF35").Value = "maa://39479 (88.89)"
$ws.Range("H37").Value = "*maa://24374 (55.56)"
$ws.Range("L37").Value = "maa://45718 (97.96), *maa://47069 (71.43), maa://45789 (100.0)"
$ws.Range("P37").Value = "maa://21280 (89.45), *maa://21239 (66.67)"
$ws.Range("AF38").Value = "maa://36697 (86.45)"
$ws.Range("H39").Value = "maa://36670 (89.11), maa://25199 (84.82), maa://30434 (91.57), ***maa://25036 (16.0), maa://45059 (83.33), *maa://44165 (66.67)"
$ws.Range("T39").Value = "maa://45788 (82.0), maa://47079 (94.12), *maa://45790 (75.0)"
$ws.Range("P41").Value = "**maa://35616 (40.0), maa://43177 (91.3)"
$ws.Range("H43").Value = "maa://22525 (92.41), maa://21284 (85.42)"
$ws.Range("H45").Value = "maa://21229 (84.82), maa://30807 (95.65), *maa://22767 (55.0), ***maa://20796 (13.79), maa://42459 (86.36)"
$ws.Range("T45").Value = "**maa://39364 (37.5)"
$ws.Range("H55").Value = "maa://32532 (91.94)"
$ws.Range("H58").Value = "*maa://37964 (60.98)"
